{"js": "// Apply the dated worksheet refresh: update the header date and every\n// two-digit-by-two-digit multiplication prompt in document order.\nconst replacements = [\n  [\"2025-09-27 Saturday\", \"2025-09-28 Sunday\"],\n  [\"83\u00d716=\", \"98\u00d793=\"],\n  [\"56\u00d734=\", \"21\u00d714=\"],\n  [\"27\u00d722=\", \"36\u00d714=\"],\n  [\"16\u00d724=\", \"39\u00d743=\"],\n  [\"52\u00d768=\", \"82\u00d718=\"],\n  [\"26\u00d787=\", \"13\u00d792=\"],\n  [\"65\u00d752=\", \"23\u00d779=\"],\n  [\"93\u00d712=\", \"59\u00d775=\"],\n  [\"13\u00d769=\", \"85\u00d733=\"],\n  [\"43\u00d744=\", \"15\u00d757=\"],\n  [\"14\u00d724=\", \"39\u00d734=\"],\n  [\"44\u00d760=\", \"94\u00d720=\"],\n  [\"41\u00d760=\", \"78\u00d736=\"],\n  [\"79\u00d725=\", \"45\u00d734=\"],\n  [\"26\u00d763=\", \"25\u00d756=\"],\n  [\"77\u00d755=\", \"76\u00d711=\"],\n  [\"85\u00d728=\", \"65\u00d711=\"],\n  [\"94\u00d763=\", \"72\u00d787=\"],\n  [\"64\u00d729=\", \"15\u00d755=\"],\n  [\"25\u00d798=\", \"30\u00d734=\"],\n  [\"22\u00d731=\", \"66\u00d761=\"],\n  [\"86\u00d760=\", \"43\u00d744=\"],\n  [\"89\u00d744=\", \"16\u00d761=\"],\n  [\"19\u00d793=\", \"43\u00d794=\"],\n  [\"41\u00d764=\", \"86\u00d783=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Each source string is unique in the document, so replace its single\n  // (first) occurrence; this also keeps behaviour correct when a later\n  // replacement's new value equals an earlier replacement's old value.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Apply the dated worksheet refresh: update the header date and every\n# two-digit-by-two-digit multiplication prompt in document order.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-09-27 Saturday\", \"2025-09-28 Sunday\"),\n  @(\"83\u00d716=\", \"98\u00d793=\"),\n  @(\"56\u00d734=\", \"21\u00d714=\"),\n  @(\"27\u00d722=\", \"36\u00d714=\"),\n  @(\"16\u00d724=\", \"39\u00d743=\"),\n  @(\"52\u00d768=\", \"82\u00d718=\"),\n  @(\"26\u00d787=\", \"13\u00d792=\"),\n  @(\"65\u00d752=\", \"23\u00d779=\"),\n  @(\"93\u00d712=\", \"59\u00d775=\"),\n  @(\"13\u00d769=\", \"85\u00d733=\"),\n  @(\"43\u00d744=\", \"15\u00d757=\"),\n  @(\"14\u00d724=\", \"39\u00d734=\"),\n  @(\"44\u00d760=\", \"94\u00d720=\"),\n  @(\"41\u00d760=\", \"78\u00d736=\"),\n  @(\"79\u00d725=\", \"45\u00d734=\"),\n  @(\"26\u00d763=\", \"25\u00d756=\"),\n  @(\"77\u00d755=\", \"76\u00d711=\"),\n  @(\"85\u00d728=\", \"65\u00d711=\"),\n  @(\"94\u00d763=\", \"72\u00d787=\"),\n  @(\"64\u00d729=\", \"15\u00d755=\"),\n  @(\"25\u00d798=\", \"30\u00d734=\"),\n  @(\"22\u00d731=\", \"66\u00d761=\"),\n  @(\"86\u00d760=\", \"43\u00d744=\"),\n  @(\"89\u00d744=\", \"16\u00d761=\"),\n  @(\"19\u00d793=\", \"43\u00d794=\"),\n  @(\"41\u00d764=\", \"86\u00d783=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  # wdFindWrap=1(wdFindContinue), wdReplace=2(wdReplaceOne) \u2014 each source\n  # string is unique in the document, so only its single occurrence changes.\n  $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
